# "include no rank decision in binary"
#
# The matrices ranking score (column G) was recomputed after removing the
# "no rank decision" cases from the binary classifier used to build the
# ranking. Column I (mat_rank, the rank position within each gender group)
# stays fixed, but because the recomputed scores changed slightly, a few
# workers now land on a different row/rank than before - so their
# prolificid/name/race/index (columns D/E/H/C) move along with them.
# All other rows keep their worker but still get an updated score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 <-> Row 3 (female rank 1 <-> rank 2) swap workers ---
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "5c5882fc5bfe7600011197cb"
$ws.Range("E2").Value = "Colleen"
$ws.Range("G2").Value = 13.42119510329043
$ws.Range("H2").Value = "White"

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E3").Value = "Annes"
$ws.Range("G3").Value = 13.17773416771519
$ws.Range("H3").Value = "Asian"

# --- Rows 4,5,6 (female rank 3,4,5) cycle workers ---
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = "608b14a312c099ac00b721b6"
$ws.Range("E4").Value = "Khushi"
$ws.Range("G4").Value = 8.277947983434146
$ws.Range("H4").Value = "Asian"

$ws.Range("C5").Value = 21
$ws.Range("D5").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("E5").Value = "Bri"
$ws.Range("G5").Value = 8.218874334828817

$ws.Range("C6").Value = 19
$ws.Range("D6").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("E6").Value = "Jewel"
$ws.Range("G6").Value = 8.21192345112825
$ws.Range("H6").Value = "Black or African American"

# --- Rows 7-13 (female rank 6-12): same workers, updated scores ---
$ws.Range("G7").Value = 5.441970684512863
$ws.Range("G8").Value = 5.381459162249058
$ws.Range("G9").Value = 5.321845954194636
$ws.Range("G10").Value = 4.498467056693604
$ws.Range("G11").Value = 4.222996349665409
$ws.Range("G12").Value = 2.390791975163696
$ws.Range("G13").Value = 1.089220531548616

# --- Rows 14-19 (male rank 1-6): same workers, updated scores ---
$ws.Range("G14").Value = 14.11239547175637
$ws.Range("G15").Value = 13.09487473480318
$ws.Range("G16").Value = 8.22111200880744
$ws.Range("G17").Value = 7.429121582096163
$ws.Range("G18").Value = 6.324528075904071
$ws.Range("G19").Value = 6.243826188088984

# --- Rows 20 <-> 21 (male rank 7 <-> 8) swap workers ---
$ws.Range("C20").Value = 32
$ws.Range("D20").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E20").Value = "Jamarii"
$ws.Range("G20").Value = 5.27722767756892
$ws.Range("H20").Value = "Black or African American"

$ws.Range("C21").Value = 33
$ws.Range("D21").Value = "60b322994d0b901954690036"
$ws.Range("E21").Value = "Brennan"
$ws.Range("G21").Value = 5.186042016282854
$ws.Range("H21").Value = "White"

# --- Rows 22-25 (male rank 9-12): same workers, updated scores ---
$ws.Range("G22").Value = 5.141087836715284
$ws.Range("G23").Value = 3.417079858592328
$ws.Range("G24").Value = 1.153463192899035
$ws.Range("G25").Value = 0.3269558257719956
